$p = $ppt.ActivePresentation

# --- 1. Refresh the "datetimeFigureOut" date placeholder text (2025-12-11 -> 2025-12-12)
#        on the slide master and every slide layout that carries one.
$oldDate = "2025-12-11"
$newDate = "2025-12-12"

$m = $p.SlideMaster
for ($i = 1; $i -le $m.Shapes.Count; $i++) {
    $shp = $m.Shapes.Item($i)
    if ($shp.HasTextFrame -and $shp.TextFrame.HasText) {
        if ($shp.TextFrame.TextRange.Text -eq $oldDate) {
            $shp.TextFrame.TextRange.Text = $newDate
        }
    }
}

for ($li = 1; $li -le $m.CustomLayouts.Count; $li++) {
    $layout = $m.CustomLayouts.Item($li)
    for ($i = 1; $i -le $layout.Shapes.Count; $i++) {
        $shp = $layout.Shapes.Item($i)
        if ($shp.HasTextFrame -and $shp.TextFrame.HasText) {
            if ($shp.TextFrame.TextRange.Text -eq $oldDate) {
                $shp.TextFrame.TextRange.Text = $newDate
            }
        }
    }
}

# --- 2. Rename the function getAll_psbl_alg_comm() <- get_all_psbl_alg_comm()
#        on every slide/shape that mentions it, without touching the sibling
#        "(...)" run that follows the identifier.
$oldName = "get_all_psbl_alg_comm"
$newName = "getAll_psbl_alg_comm"

for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $slide = $p.Slides.Item($si)
    for ($shi = 1; $shi -le $slide.Shapes.Count; $shi++) {
        $shp = $slide.Shapes.Item($shi)
        if ($shp.HasTextFrame -and $shp.TextFrame.HasText) {
            $tr = $shp.TextFrame.TextRange
            $text = $tr.Text
            $idx = $text.IndexOf($oldName)
            if ($idx -ge 0) {
                $chars = $tr.Characters($idx + 1, $oldName.Length)
                $chars.Text = $newName
            }
        }
    }
}
